# Update cryptocurrency price/volume data per upstream refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "95.286.07"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.85%  "
$ws.Range("E2").Style = "Normal"

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.551.63"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.65%  "
$ws.Range("E3").Style = "Normal"

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("E4").Style = "Normal"

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.40"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.68%  "
$ws.Range("E5").Style = "Normal"

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "648.40"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +2.62%  "
$ws.Range("E6").Style = "Normal"

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.45"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +1.22%  "
$ws.Range("E7").Style = "Normal"

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.397"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("E8").Style = "Normal"

# Row 9
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("E9").Style = "Normal"

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.986"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.82%  "
$ws.Range("E10").Style = "Normal"

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "3.547.04"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.61%  "
$ws.Range("E11").Style = "Normal"

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.201"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.35%  "
$ws.Range("E12").Style = "Normal"

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "41.93"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -3.70%  "
$ws.Range("E13").Style = "Normal"

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.51"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +3.88%  "
$ws.Range("E14").Style = "Normal"

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.227.96"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.91%  "
$ws.Range("E15").Style = "Normal"

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "95.173.59"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.91%  "
$ws.Range("E16").Style = "Normal"

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000251"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.21%  "
$ws.Range("E17").Style = "Normal"

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.556.56"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.59%  "
$ws.Range("E18").Style = "Normal"

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.87"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -4.96%  "
$ws.Range("E19").Style = "Normal"

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.68"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.05%  "
$ws.Range("E20").Style = "Normal"

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.68"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.83%  "
$ws.Range("E21").Style = "Normal"

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.47"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +2.59%  "
$ws.Range("E22").Style = "Normal"

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "503.33"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -2.09%  "
$ws.Range("E23").Style = "Normal"

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.471"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -4.55%  "
$ws.Range("E24").Style = "Normal"

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000192"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.16%  "
$ws.Range("E25").Style = "Normal"

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.54"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.14%  "
$ws.Range("E26").Style = "Normal"

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "94.83"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +2.56%  "
$ws.Range("E27").Style = "Normal"

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.42"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.72%  "
$ws.Range("E28").Style = "Normal"

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.749.26"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.84%  "
$ws.Range("E29").Style = "Normal"

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.04"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.34%  "
$ws.Range("E30").Style = "Normal"

# Row 31
$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = "Dai"
$ws.Range("B31").Style = "Normal"
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("C31").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("E31").Style = "Normal"

# Row 32
$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("B32").Style = "Normal"
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("C32").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.19"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.89%  "
$ws.Range("E32").Style = "Normal"

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.138"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -3.52%  "
$ws.Range("E33").Style = "Normal"

# Row 34
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.32%  "
$ws.Range("E34").Style = "Normal"

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.177"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -2.53%  "
$ws.Range("E35").Style = "Normal"

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "31.83"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +5.31%  "
$ws.Range("E36").Style = "Normal"

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.556"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.62%  "
$ws.Range("E37").Style = "Normal"

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.14"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +7.51%  "
$ws.Range("E38").Style = "Normal"

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "560.62"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -3.85%  "
$ws.Range("E39").Style = "Normal"

# Row 40
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +4.19%  "
$ws.Range("E40").Style = "Normal"

# Row 41
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.00%  "
$ws.Range("E41").Style = "Normal"

# Row 42
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.51%  "
$ws.Range("E42").Style = "Normal"

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.898"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -3.41%  "
$ws.Range("E43").Style = "Normal"

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.76"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.69%  "
$ws.Range("E44").Style = "Normal"

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "34.95"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +34.35%  "
$ws.Range("E45").Style = "Normal"

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.30"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +6.37%  "
$ws.Range("E46").Style = "Normal"

# Row 47
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "WhiteBITCoin"
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "23.55"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.90%  "
$ws.Range("E47").Style = "Normal"

# Row 48
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "Filecoin"
$ws.Range("B48").Style = "Normal"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.61"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.18%  "
$ws.Range("E48").Style = "Normal"

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0410"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -3.13%  "
$ws.Range("E49").Style = "Normal"

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.54"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.27%  "
$ws.Range("E50").Style = "Normal"

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "53.34"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.03%  "
$ws.Range("E51").Style = "Normal"

